$d = $word.ActiveDocument

# --- 1. Replace the page-break paragraph (right after 'June 2024') with two empty paragraphs ---
$breakPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $t = $pp.Range.Text
    if ($t.Length -eq 2 -and [int][char]$t[0] -eq 12 -and [int][char]$t[1] -eq 13) {
        $breakPara = $pp
        break
    }
}
if ($breakPara -ne $null) {
    $bStart = $breakPara.Range.Start
    $bEnd = $breakPara.Range.End
    $r = $d.Range($bStart, $bEnd)
    $r.Text = "`r`r"
}

# --- 2. Replace the four body paragraphs with their revised text ---
$newIntro = @'
This study explores the intersection of gamification, personalization, and adaptivelearning. The study uses the Decision Tree Thompson Sampling (DTTS) algorithm into gamified content. It aims to evaluate the algorithm's effectiveness in personalized learning, and exploring its adaptability across diverse subjects. The research aims to drive positive transformations in educational practices, ultimatelybenefiting both educators and students in a rapidly changing educational landscape. This study is a precedent for future research on theintegration of advanced educational algorithms in educational settings. This chapter outlines the background, problem statement and significance of the study. The study is purposefully limited in scope to ensure clarity and coherence in its objectives. The research on gamification in education, specifically employing the Decision Tree Thompson Sampling (DTTS) algorithm, holds significance in several key dimensions. This approach helps in transforming traditional teaching methodologies, increasing student participation and interest. The shift towards a more engaging, personalized, and inclusive educational environment marks a positive step forward in overcoming the complex challenges faced by educational institutions and educators.
'@
$newMethod = @'
The study aims to assess the effectiveness and adaptability of the ALGE-BRUH webgame in providing personalized educational content. The main objective of the study is to modify the online educational game to better suit Grade 10 students' unique learning needs and preferences. The game uses the Multi-Armed Bandit (MAB) and the DTTS (Dynamic Time-Triggered Switching) Algorithm to enhance learning outcomes and student engagement. The data collection approach for this study will be divided into two phases, each purposefully intended to capture both user information and evaluations comprehensively. The ALGE-BRUH webgame leverages real-time user data to inform the DTTS algorithm, which personalizes the learning experience. The MAB Algorithm optimizes the selection of learningcontent and activities based on students' performance metrics. The study focuses on evaluating the algorithm's effectiveness in personalized learning, understanding student satisfaction, and enjoyment. The research will move into the second phase, which involves collecting user assessments via carefully constructed questionnaires, after user data has been collected. Gamification refers to the use of game-like elements, such as rewards, progress levels, and interactive structures, in educational settings. The Multi-Armed Bandit Algorithm incorporates a heuristic method that focuses on how users have answered questionscorrectly and quickly. The selection of the new system's concept will be predicated upon the insights derived from this research endeavor. The project will be broken down into smaller tasks. This strategy ensures continuous improvement and effective adaptation within the gamified educational environment.
'@
$newResults = @'
The ALGE-BRUH webgame is designed to enhance educational outcomes and support academic achievement in Grade 10 algebra. The system offers real-time analytics capabilities to monitor and analyze student performance and progress over time. The game is designed with HTML5, CSS, and JavaScript. It is powered by the DTTS algorithm. The project structure contains screenshots and its description with its description in its description. The software is currently in development and is expected to be completed by the end of the year. The project involves the development of a web-based game system named ALGE-BRUH. It is designed to integrate the Decision Tree Thompson Sampling (DTTS) algorithm into gamified educational content. The system is built using HTML, CSS, and JavaScript for the front-end interface and game logic. The majority of students reported that the game helped improve their problem-solving abilities and overall academic performance. The design approach aligns with popular internet culture to make the game engaging for Grade 10 page and Sign in page users. The majority of students reported positive experiences, with high percentages indicating they felt content, skilled, engaged, and found the game enjoyable. The engagement and enjoyment from the game likely contribute to higher motivation and sustained interest in learning. The system’s ability to offer real-time analytics further enhances its performance. It allows continuous monitoring and adaptive responses based on user interactions. The user can use items to give the current character a ‘power up’. The player can use the ‘Power Up’ button to give their character a power up.
'@
$newDiscussion = @'
60% of the students strongly agreed that the game could be valuable for their learning. 62.22% strongly agreed the game was very interesting. 53.33% felt skilled at playing the game, although 22.22%, felt they could not do very well. The Multi-Armed Bandit algorithm was effective in creating an engaging, enjoyable, and educational experience for the students. The game was perceived as valuable and useful for learning, particularly in enhancing math skills, problem-solving abilities, and overall academic performance. 55.56% felt competent after playing the game for a while. to other students, and 55.56 per cent felt competent to play the game with other students. The game was designed to test students' knowledge of the game. It was designed for students to test their knowledge of a game and its rules. It is not intended to be a test of knowledge of game rules. The test was not meant to be an assessment of the player's knowledge of these rules, but of the games rules.
'@

$introOld = "This research is driven"
$methodOld = "The study aims to assess the effectiveness and adaptability of the ALGE-BRUH webgame"
$resultsOld = "The ALGE-BRUH webgame demonstrates strong performance"
$discussionOld = "60% of the students strongly agreed that the game could be valuable"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $t = $pp.Range.Text
    if ($t.StartsWith($introOld)) { $pp.Range.Text = $newIntro }
    elseif ($t.StartsWith($methodOld)) { $pp.Range.Text = $newMethod }
    elseif ($t.StartsWith($resultsOld)) { $pp.Range.Text = $newResults }
    elseif ($t.StartsWith($discussionOld)) { $pp.Range.Text = $newDiscussion }
}

Write-Output "done"
